# Update InboxData.xlsx ("FunctionAndInbox" sheet) test data.
# Old layout:  A1=Function   B1=Inbox
#              A2=Inbound    B2=Purchase Orders
# New layout:  A1=Persona    B1=Inbox        C1=SearchRecord
#              A2=Customer Success  B2=Sales Orders  C2=ZOR
#
# Cell-write order below intentionally matches the order new shared-string
# values first appear in the target workbook (Inbox, Persona,
# Customer Success, Sales Orders, SearchRecord, ZOR) so the rebuilt
# sharedStrings table lines up with the authored diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = "Inbox"
$ws.Range("A1").Value = "Persona"
$ws.Range("A2").Value = "Customer Success"
$ws.Range("B2").Value = "Sales Orders"
$ws.Range("C1").Value = "SearchRecord"
$ws.Range("C2").Value = "ZOR"

# Column A widens to fit the new "Customer Success" / "Persona" text.
$null = $ws.Range("A1").EntireColumn.AutoFit()

# New last column becomes the active selection, matching the saved view.
$null = $ws.Range("C1").Select()
